# Finalização do relatório final
# Updates the "Planilha1" rule-name / metadata text to corrected/expanded
# wording, widens column A, moves the active selection, and (as a
# consequence of the text changes) the VLOOKUP-driven "tabela resumo"
# sheet now returns #N/A for the rows whose source text changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- sheet view: drop the frozen topLeftCell, move the active selection ---
$ws.Range("A4").Select()

# --- widen column A (raw stored width 74 == ColumnWidth 74 - 5/6) ---
$ws.Columns.Item(1).ColumnWidth = 73.16666666666667

# --- row 4 grows from 2 lines to 3 lines tall ---
$ws.Rows.Item(4).RowHeight = 45
# --- rows 13 & 14 gain an explicit height matching the rest of the table ---
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30

# --- column B (Metadado) text corrections: lowercase field names capitalised ---
$ws.Range("B2").Value = ",Measurements, Measurements_Altura, Measurements_Largura, Measurements_Profundidade, Measurements_Diametro, Measurements_Peso, Measurements_Espessura,"
$ws.Range("B5").Value = ",Class,Creation Location,Creator,Description,Location,Materials and Techniques,Other Descriptive Notes,Title,Work Type,"
$ws.Range("B6").Value = ",Class,Creation Location,Creator,Inscription,Location,Materials and Techniques,Measurements, Measurements_Altura, Measurements_Largura, Measurements_Profundidade, Measurements_Diametro, Measurements_Peso, Measurements_Espessura,Physical Description,Work Type,Inscription,Location,"
$ws.Range("B7").Value = ",Measurements, Measurements_Altura, Measurements_Largura, Measurements_Profundidade, Measurements_Diametro, Measurements_Peso, Measurements_Espessura,"
$ws.Range("B8").Value = ",Creator,Inscription,Materials and Techniques,Measurements, Measurements_Altura, Measurements_Largura, Measurements_Profundidade, Measurements_Diametro, Measurements_Peso,Measurements_Espessura,Work Type,Title,Date,Location,"
$ws.Range("B9").Value = ",Measurements, Measurements_Altura, Measurements_Largura, Measurements_Profundidade, Measurements_Diametro, Measurements_Peso, Measurements_Espessura,"
$ws.Range("B16").Value = ",Class,Materials and Techniques,Work Type,"
$ws.Range("B18").Value = ",Measurements, Measurements_Altura, Measurements_Largura, Measurements_Profundidade, Measurements_Diametro, Measurements_Peso, Measurements_Espessura,"

# --- column A (Nome da Regra) text corrections ---
$ws.Range("A3").Value = "Anos com menos de 4 digitos, inserir 0 à esquerda"
$ws.Range("A4").Value = "Capitalizar as inicais de nomes próprios e da primeira palavra, para outros termos use letras minúsculas"
$ws.Range("A13").Value = "Seguir padrão para registro de horas, minutos e segundos"
$ws.Range("A14").Value = "Seguir padrão para registro de dia, mês e ano da data"

# --- "tabela resumo": the VLOOKUPs for the rows whose Planilha1 rule-name
# text changed (old text no longer found) now resolve to #N/A automatically
# once the lookup source text above changes; the rows merely need their
# formulas to recalculate (handled automatically by the engine). Two of the
# sheet's rows also lose their explicit row height override, reverting to
# the sheet's natural/auto height. ---
$ws2 = $wb.Worksheets.Item("tabela resumo")
$ws2.Rows.Item(11).AutoFit()
$ws2.Rows.Item(16).AutoFit()
